# Auto-generated Excel COM-interop edit script
# Commit: Extension Installation Green Tech from EU to all the regions

$wb = $excel.ActiveWorkbook
$wsAnnual = $wb.Worksheets.Item("Annual production")
$wsCumulative = $wb.Worksheets.Item("Cumulative")

# --- Rename the "Copper ores and concentrates" label to "Copper" ---
# (update on both sheets so the shared string text is replaced everywhere)
$wsAnnual.Range("C4").Value = "Copper"
$wsCumulative.Range("C4").Value = "Copper"

# --- New annual-production values for rows 2-5 (columns D:CO, years 2011-2100) ---
$newRow2 = @(1072.373874491149, 1397.73011970385, 1397.031627798344, 1327.590969412265, 1780.597773721523, 1284.725912258199, 2324.625080083905, 2311.609606929499, 2058.411181797298, 1901.455711851288, 4070.010987434653, 2490.436145986909, 4975.54675684899, 5356.238079114935, 5726.937853245291, 12512.08236678371, 13425.38294233107, 14354.25906959043, 15310.38161680337, 16285.53220692, 16397.08394224583, 17320.12292895656, 18258.56367421243, 19166.01462993056, 20123.35243165829, 9663.057042433004, 9991.49481271195, 10374.09756956074, 10639.67017657794, 10975.96016470608, 7514.689765459631, 8133.71366963574, 8748.139778199013, 9800.763977031189, 11204.74498363686, 15053.15457487533, 17518.08986132502, 19060.05398536073, 20045.17974835873, 20909.08286616235, 24903.16018992295, 25600.25593945458, 26300.82205184953, 26660.6553759644, 26076.6911099267, 24686.74149987378, 22017.97669553811, 21200.67582481223, 21307.18519446514, 21160.76913727244, 19693.38985203286, 19132.43988275844, 19503.40766641527, 20583.45892257623, 22194.66141449579, 23758.39924256544, 26250.97185823477, 28308.65310447558, 29867.67819708808, 31348.07298804109, 32733.55852676531, 34493.91599636365, 35505.535283206, 35883.81691142748, 35691.72625755112, 34627.67682837267, 33297.54749183558, 32107.32794683842, 31561.6070794988, 31331.25254896398, 30717.26757958199, 30311.16249037089, 30442.19335836857, 31269.30307931092, 32630.56462613314, 34132.00870267345, 36053.97470648385, 38054.55202523863, 39861.01736833101, 41446.69170639732, 42675.74019383219, 44070.13912861065, 45206.06613414158, 45805.24189583262, 45794.70966571696, 44834.91029193866, 43895.15184659533, 42917.90311027105, 42184.0380882688, 41748.54603022619)
$newRow3 = @(70.91124739014649, 92.70500158204401, 92.53982554733052, 87.87484380016343, 118.398417739643, 85.08204867090255, 154.8130776054375, 153.7869514538161, 136.8280362575191, 126.2151906065665, 271.7826109078967, 165.8299747707763, 332.6222416492287, 358.1634799391416, 383.0428350747089, 838.6005569577648, 899.9248120558132, 962.2987289860444, 1026.522512639116, 1092.027396853858, 1099.49883781615, 1161.51144535633, 1224.579344083275, 1285.628267841386, 1350.037133737828, 647.5721933183137, 669.526783918712, 695.2682577772957, 713.210029492169, 735.752032137933, 503.4660975158289, 545.2712455528898, 586.9684537731002, 658.3459520437536, 753.493389753566, 1013.277468459089, 1180.445852067391, 1284.912455596786, 1351.533341074027, 1409.955409286491, 1678.551565381215, 1725.817953132386, 1774.056384195977, 1799.949479940971, 1759.482530182979, 1663.171188687445, 1480.552439428979, 1424.360174647541, 1431.149893651275, 1420.855470174374, 1321.813330130237, 1283.833046019663, 1309.108050681848, 1382.590695788589, 1492.186565328819, 1599.028100735882, 1768.475376335078, 1908.391340662781, 2014.496853834482, 2115.725203871159, 2211.350799830325, 2331.327429565635, 2398.85816641544, 2422.95938669616, 2408.154928489769, 2334.07409814318, 2241.885882562733, 2159.808659714262, 2122.156109580928, 2106.339754630997, 2064.862044798705, 2037.452030461959, 2046.60669034907, 2103.187279821532, 2196.281895947504, 2299.396487131841, 2431.448595860956, 2569.326249763149, 2693.54604194142, 2801.21414216088, 2883.753401727607, 2976.900893939906, 3052.161771546927, 3090.72540504576, 3087.985073013786, 3021.466675259153, 2956.431608188964, 2889.212592497669, 2838.892662785518, 2809.195503291899)
$newRow4 = @(1624231.632163128, 1459073.662060171, 1273950.194004356, 1477943.369022302, 1624319.753844372, 1544497.383707435, 1680352.060849322, 1544436.16825248, 1539632.871095464, 2062866.995901393, 1982943.027974252, 1977433.52891767, 9988456.547842104, 10023206.77175408, 10067016.18517715, 31346381.60078163, 31434058.19355349, 31530328.61320518, 31638593.91209295, 31751126.67332055, 32705068.83671846, 32824733.992457, 32954563.23460107, 33107401.98981949, 33321025.76875351, 9468884.231811421, 9673112.529013738, 9959024.689114314, 10363500.40932948, 10903771.98552297, 3537279.823616583, 4218116.132766244, 5132803.330663079, 6336331.69868101, 7658758.673808551, 12575153.28875133, 14236851.86783882, 15392268.56234481, 16197664.71244349, 16824664.33864639, 22245086.17041712, 22426917.33226267, 22432491.35557772, 22170577.71436144, 21480054.51213252, 21586560.45610342, 20148566.23603553, 19247218.69660724, 18669042.86510228, 18136063.16725055, 16693967.08568335, 16358435.41742795, 16411376.84401968, 16859876.66356125, 17628774.38081717, 18136497.5770331, 19337936.29040986, 20426287.16384293, 21313566.05727978, 22104533.8900356, 22769378.62236677, 23485259.98516439, 23938970.27027962, 24147820.17159694, 24106791.16916629, 23767073.18623966, 23255623.22083038, 22760897.70991759, 22513442.74228181, 22351448.58683594, 22020827.35374447, 21829735.87213953, 21731852.73189611, 21878360.32475173, 22313559.56553872, 22690901.35739988, 23384370.55862211, 24093816.05285225, 24720963.16372762, 25249204.11217924, 25362114.75606243, 25799068.67806806, 26152469.96482579, 26355037.88191808, 26385560.31508645, 25688302.89111708, 25445898.45172931, 25204074.63607221, 25240580.48175715, 25263766.67022805)
$newRow5 = @(2150342.287831458, 1440641.974271039, 1347505.296692689, 1300054.160795987, 1248697.590957438, 1561555.338495209, 1721155.079761696, 1504823.440738913, 1399739.683168624, 1573380.696318417, 1665942.493873482, 2164842.236706774, 16378591.24583179, 16424048.39200231, 16479178.63927891, 54346972.91804177, 54424605.85412651, 54515923.13989751, 54624375.55369617, 54756060.90755674, 57718555.44862213, 57930912.93346493, 58211295.78807408, 58586852.64134628, 59092060.36553384, 18050986.226274, 18947156.92232763, 20114991.27481109, 21606417.32356417, 23467516.09420415, 11302349.00873467, 13983163.3308182, 17066347.62622076, 20503333.49128058, 24206747.71491146, 33954451.65901872, 37775487.39658225, 41388388.15400425, 44598940.75017184, 47224879.10988516, 57320167.46640705, 58381140.0911729, 58577958.01154337, 57946830.97043698, 56589391.12816014, 57061084.73049114, 54759761.43091525, 52303284.38718721, 49916091.9942808, 47811530.14230615, 44778584.81409178, 43753678.4235727, 43427570.76142555, 43826311.95617687, 44915751.33440275, 45839232.00281735, 48015226.56886657, 50525984.75771311, 53215275.3615213, 55928484.01214515, 58416492.25641885, 60765287.81333553, 62770791.07290171, 64365754.53698423, 65517178.82686125, 66250484.04016297, 66550333.71226293, 66501180.70898025, 66184742.02709807, 65697032.51440839, 64920881.28671881, 64395968.65340099, 63993625.19470124, 63788705.24571463, 63835375.79006504, 63773535.86794692, 64391802.80947018, 65283129.74672952, 66411358.50234376, 67724596.04027189, 68490091.485071, 69981934.47659098, 71463779.65164584, 72875819.15731166, 74168271.01054494, 74251876.93665504, 75208976.01531665, 75979735.17540646, 76570981.81259938, 77002064.43922363)

function Set-RowValues($ws, [int]$rowNum, $values) {
    $n = $values.Length
    $arr = New-Object "object[,]" 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $startCell = $ws.Cells.Item($rowNum, 4)   # column D
    $endCell   = $ws.Cells.Item($rowNum, 3 + $n)
    $ws.Range($startCell, $endCell).Value = $arr
    return $values
}

function Set-CumulativeRowValues($ws, [int]$rowNum, $values) {
    $n = $values.Length
    $arr = New-Object "object[,]" 1,$n
    $running = 0.0
    for ($i = 0; $i -lt $n; $i++) {
        $running += [double]$values[$i]
        $arr[0,$i] = $running
    }
    $startCell = $ws.Cells.Item($rowNum, 4)   # column D
    $endCell   = $ws.Cells.Item($rowNum, 3 + $n)
    $ws.Range($startCell, $endCell).Value = $arr
}

# Apply to "Annual production" sheet (rows 2-5 map to Neodymium, Dysprosium, Copper, Raw silicon)
Set-RowValues $wsAnnual 2 $newRow2 | Out-Null
Set-RowValues $wsAnnual 3 $newRow3 | Out-Null
Set-RowValues $wsAnnual 4 $newRow4 | Out-Null
Set-RowValues $wsAnnual 5 $newRow5 | Out-Null

# Recompute "Cumulative" sheet as the row-wise running total of the Annual production sheet
Set-CumulativeRowValues $wsCumulative 2 $newRow2
Set-CumulativeRowValues $wsCumulative 3 $newRow3
Set-CumulativeRowValues $wsCumulative 4 $newRow4
Set-CumulativeRowValues $wsCumulative 5 $newRow5

